# Docs: Update Diagrams and PPTX replace browser panel with person panel
#
# The UiComponentClassDiagram slide has a small rectangle labelled
# "BrowserPanel" (component box for the old browser UI panel). The UI was
# renamed, so the box now represents the "PersonPanel" component instead.
# Find that shape by its current text (robust against z-order/id churn)
# and rename it in place, leaving its formatting untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$oldText = "BrowserPanel"
$newText = "PersonPanel"

$shapeCount = $s.Shapes.Count
for ($i = 1; $i -le $shapeCount; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq $oldText) {
            $tr.Text = $newText
        }
    }
}
